$d = $word.ActiveDocument
$bullet = [char]0x2022

# ------------------------------------------------------------------
# 1) CORE COMPETENCIES: collapse the three detailed bullet paragraphs
#    into a single condensed paragraph.
# ------------------------------------------------------------------
$firstIdx = 0
$secondIdx = 0
$thirdIdx = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($firstIdx -eq 0 -and $t.StartsWith("Statistical Analysis & Machine Learning:")) {
        $firstIdx = $i
    }
    elseif ($secondIdx -eq 0 -and $t.StartsWith("Big Data & Data Engineering:")) {
        $secondIdx = $i
    }
    elseif ($thirdIdx -eq 0 -and $t.StartsWith("Data Visualization & Reporting:")) {
        $thirdIdx = $i
    }
}

if ($firstIdx -gt 0 -and $secondIdx -gt 0 -and $thirdIdx -gt 0) {
    # Delete the 2nd and 3rd paragraphs entirely (including their paragraph marks).
    $delRange = $d.Range($d.Paragraphs($secondIdx).Range.Start, $d.Paragraphs($thirdIdx).Range.End)
    $delRange.Delete()

    # Replace the remaining (first) paragraph's text with the condensed summary.
    $condensed = "Statistical Analysis & Machine Learning $bullet Big Data & Data Engineering $bullet Data Visualization & Reporting"
    $d.Paragraphs($firstIdx).Range.Text = $condensed
}

# ------------------------------------------------------------------
# 2) Append a new "TECHNICAL SKILLS" section at the very end of the
#    document body (before the final section break).
# ------------------------------------------------------------------
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)

$line1 = "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques"
$line2 = "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization"
$line3 = "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation"

$r.InsertAfter("`rTECHNICAL SKILLS`r$line1`r$line2`r$line3")

$newCount = $d.Paragraphs.Count
$headingPara = $d.Paragraphs($newCount - 3)
$headingPara.Style = "Heading2"

Write-Output "OK"
